# "added for User Module"
# Update the "User" worksheet (sheet3) with new/changed test data rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("User")
$ws.Activate()

# Existing endpoint id values change from "U37" to new, distinct ids.
$ws.Range("B6").Value = "U52"
$ws.Range("B4").Value = "U60"

# Row 6 gains two more "Active" status columns (P6, Q6).
$ws.Range("P6").Value = "Active"
$ws.Range("Q6").Value = "Active"

# New row 8: a "Delete" scenario entry with its own user id.
$ws.Range("A8").Value = "Delete"
$ws.Range("B8").Value = "U141"

# Leave the selection on the newly added cell, matching the saved file.
$ws.Range("B8").Select()
